# Edit script: 
#  1) Re-style the three tables (slides 14, 15, 16) from the old built-in
#     table style GUID to the new one.
#  2) Swap the presentation's theme color scheme from "Integral" (Red
#     Violet) to the standard "Office Theme" colors.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{F7C970EF-2A9C-4AE5-B220-296AD76D5FF6}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
